$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing MSSQL rows (5-8, which are sheet rows 6-8 plus row5=Oracle unaffected)
# Row 6 (Database/MSSQL): rename to shorter version string
$ws.Range("C6").Value = "Microsoft SQL Server 2019"
$ws.Range("C7").Value = "Microsoft SQL Server 2017"
$ws.Range("C8").Value = "Microsoft SQL Server 2016"

# Add new rows: MSSQL 2022, and two MariaDB rows
$ws.Range("A9").Value = "Database"
$ws.Range("B9").Value = "MSSQL"
$ws.Range("C9").Value = "Microsoft SQL Server 2022"
$ws.Range("D9").Value = "CIS"

$ws.Range("A10").Value = "Database"
$ws.Range("B10").Value = "MARIA"
$ws.Range("C10").Value = "MariaDB 10_11"
$ws.Range("D10").Value = "CIS"

$ws.Range("A11").Value = "Database"
$ws.Range("B11").Value = "MARIA"
$ws.Range("C11").Value = "MariaDB 10_6"
$ws.Range("D11").Value = "CIS"

# Update selection to match final state (active cell C11)
$ws.Range("C11").Select()
